# The presentation ships two DrawingML theme parts:
#   ppt/theme/theme1.xml  -> bound to the (single) slide master, "Integral" /
#                            "Red Violet" colour scheme
#   ppt/theme/theme2.xml  -> bound to the notes master, default "Office Theme"
#                            colour scheme
#
# The target revision swaps the two themes' content (theme1 becomes the
# Office Theme colours, theme2 becomes the Integral/Red Violet colours).
# Font scheme and format scheme are identical between the two themes, so the
# only real content difference is the 12-slot colour scheme (dk1, lt1, dk2,
# lt2, accent1-6, hlink, folHlink).
#
# Apply the reachable half of that swap through the PowerPoint object model:
# recolour the deck's theme (theme1.xml, via the slide master / any slide's
# ThemeColorScheme) to the Office Theme palette.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# index -> (role, target RGB as 0xBBGGRR for the COM RGB property)
$tcs.Colors(1).RGB  = 0x000000   # dk1      -> 000000
$tcs.Colors(2).RGB  = 0xFFFFFF   # lt1      -> FFFFFF
$tcs.Colors(3).RGB  = 0x6A5444   # dk2      -> 44546A
$tcs.Colors(4).RGB  = 0xE6E6E7   # lt2      -> E7E6E6
$tcs.Colors(5).RGB  = 0xD59B5B   # accent1  -> 5B9BD5
$tcs.Colors(6).RGB  = 0x317DED   # accent2  -> ED7D31
$tcs.Colors(7).RGB  = 0xA5A5A5   # accent3  -> A5A5A5
$tcs.Colors(8).RGB  = 0x00C0FF   # accent4  -> FFC000
$tcs.Colors(9).RGB  = 0xC47244   # accent5  -> 4472C4
$tcs.Colors(10).RGB = 0x47AD70   # accent6  -> 70AD47
$tcs.Colors(11).RGB = 0xC16305   # hlink    -> 0563C1
$tcs.Colors(12).RGB = 0x724F95   # folHlink -> 954F72
